# Add a new species ("European smelt" / "Osmerus eperlanus") to both the
# FR (sheet1 / Tableau1) and EN (sheet2 / Tableau13) tables, by cloning the
# 20-row "Phase x Sex" block of an existing species and swapping the
# species name (col A) / Latin name (col E).

$wb = $excel.ActiveWorkbook
$wsFR = $wb.Worksheets.Item(1)
$wsEN = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1 / FR: clone rows 822:841 (Grey gurnard / Eutrigla gurnardus) down to
# 842:861, then relabel col A/E to the new species.
# ---------------------------------------------------------------------------
$wsFR.Range("A822:E841").Copy($wsFR.Range("A842"))
$wsFR.Range("A842:A861").Value = "European smelt"
$wsFR.Range("E842:E861").Value = "Osmerus eperlanus"

# Column E on the new FR rows uses the bordered style (same as E50, say)
# rather than the plain style the row-clone produced.
$wsFR.Range("E50").Copy()
$wsFR.Range("E842:E861").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$loFR = $wsFR.ListObjects.Item(1)
$loFR.Resize($wsFR.Range("A1:E861"))

# ---------------------------------------------------------------------------
# Sheet2 / EN: clone rows 816:835 (Gilt-head bream / Sparus aurata) down to
# 836:855, then relabel col A/E to the new species.
# ---------------------------------------------------------------------------
$wsEN.Range("A816:E835").Copy($wsEN.Range("A836"))
$wsEN.Range("A836:A855").Value = "European smelt"
$wsEN.Range("E836:E855").Value = "Osmerus eperlanus"

$loEN = $wsEN.ListObjects.Item(1)
$loEN.Resize($wsEN.Range("A1:E855"))

# ---------------------------------------------------------------------------
# View state: EN becomes the active tab, with new selections/scroll
# positions on both sheets mirroring the appended rows.
# ---------------------------------------------------------------------------
$wsFR.Activate()
$excel.ActiveWindow.ScrollRow = 849
$excel.ActiveWindow.ScrollColumn = 1
$wsFR.Range("A842:A861").Select()

$wsEN.Activate()
$excel.ActiveWindow.ScrollRow = 832
$excel.ActiveWindow.ScrollColumn = 1
$wsEN.Range("D847").Select()
